$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-25.
$data = @{
    2  = @(5, 5)
    3  = @(7, 8)
    4  = @(7, 8)
    5  = @(7, 9)
    6  = @(4, 6)
    7  = @(7, 8)
    8  = @(1, 4)
    9  = @(1, 5)
    10 = @(1, 4)
    11 = @(1, 3)
    12 = @(1, 6)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(2, 6)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(6, 7)
    23 = @(3, 5)
    24 = @(1, 3)
    25 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
